# "add genders in batumi" — add a new 2023 data column (S) to the yearly
# trade table, and a fresh trailing blank column (U) that mirrors the
# previous blank "spacer" column (T), shifting the used range out by one
# column (T39 -> U39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- New 2023 values for rows 4-14 (row 3 is the year header itself) ---
$values = @{
    4  = 1132.8
    5  = 182.9
    6  = 6970
    7  = 4164
    8  = 733.5
    9  = 36
    10 = 37.1
    11 = 147
    12 = 7.9
    13 = 999.1
    14 = 965
}

# Row 3: year header "2023" in S3, formatted like R3 (the 2022 header cell).
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial($xlPasteFormats)
$ws.Range("S3").Value = 2023

# Rows 4-14: new data values in column S, formatted like the existing R column.
foreach ($row in $values.Keys) {
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial($xlPasteFormats)
    $ws.Range("S$row").Value = $values[$row]
}

# Row 1 (merged title band) grows by one cell: format it like its neighbour.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial($xlPasteFormats)

# Column U becomes the new trailing blank spacer column, mirroring the
# formatting that column T already had in every row of the table.
$ws.Range("T2").Copy()
$ws.Range("U2").PasteSpecial($xlPasteFormats)

for ($row = 3; $row -le 15; $row++) {
    $ws.Range("T$row").Copy()
    $ws.Range("U$row").PasteSpecial($xlPasteFormats)
}

$ws.Application.CutCopyMode = $false

# Extend the header merge A1:R1 -> A1:S1.
$ws.Range("A1:S1").Merge()

# Match the selection left behind by the edit (2023 column highlighted).
$ws.Range("S3:S14").Select()
